# Auto-generated edit script to update Trials and Summary sheets
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trials")
$wsSummary = $wb.Worksheets.Item("Summary")

# Each entry: row, Iterations(B), Runtime(C), Success(D or $null if unchanged)
$updates = @(
    @(2, 4420, 2.767275333404541, $null),
    @(3, 3027, 1.62413215637207, $null),
    @(4, 11046, 12.73072671890259, 1),
    @(5, 3088, 1.765993118286133, 1),
    @(6, 122, 0.05011892318725586, $null),
    @(7, 9659, 11.12273788452148, $null),
    @(8, 3426, 2.291496515274048, $null),
    @(9, 11418, 14.83373045921326, $null),
    @(10, 6575, 6.192307472229004, $null),
    @(11, 4727, 4.142014503479004, $null),
    @(12, 1362, 0.7508976459503174, $null),
    @(13, 1928, 1.178695440292358, $null),
    @(14, 3210, 1.52738618850708, $null),
    @(15, 2621, 1.427172899246216, $null),
    @(16, 6833, 6.65604043006897, 1),
    @(17, 11951, 15.51112461090088, $null),
    @(18, 4889, 3.699402570724487, $null),
    @(19, 10767, 15.81407856941223, $null),
    @(20, 13763, 21.3364269733429, $null),
    @(21, 4023, 2.77872109413147, $null),
    @(22, 8586, 8.907962799072266, $null),
    @(23, 1815, 1.02644157409668, $null),
    @(24, 16000, 28.41391444206238, 1),
    @(25, 25507, 60.00278210639954, 0),
    @(26, 4393, 2.884665012359619, $null),
    @(27, 3584, 2.055037260055542, $null),
    @(28, 21118, 40.36164212226868, 1),
    @(29, 3183, 1.670277833938599, $null),
    @(30, 6493, 4.666834831237793, $null),
    @(31, 7765, 6.401481628417969, 1),
    @(32, 3200, 1.806713104248047, $null),
    @(33, 11730, 13.1047625541687, $null),
    @(34, 2491, 1.294078826904297, $null),
    @(35, 15405, 24.74996066093445, $null),
    @(36, 1901, 0.8401088714599609, $null),
    @(37, 14090, 19.17671704292297, $null),
    @(38, 5426, 4.124681949615479, $null),
    @(39, 7105, 5.887673616409302, $null),
    @(40, 9311, 10.53630137443542, $null),
    @(41, 787, 0.2993359565734863, $null),
    @(42, 5798, 4.166477680206299, $null),
    @(43, 16174, 28.6018385887146, $null),
    @(44, 10587, 12.39374470710754, 1),
    @(45, 3404, 2.172279596328735, $null),
    @(46, 2469, 1.071197509765625, $null),
    @(47, 10953, 12.21662950515747, $null),
    @(48, 4264, 2.997926235198975, $null),
    @(49, 7987, 8.435249090194702, $null),
    @(50, 8553, 8.555781841278076, $null),
    @(51, 1150, 0.5058367252349854, $null),
    @(52, 1718, 0.8110020160675049, 1),
    @(53, 20764, 44.52081203460693, $null),
    @(54, 9624, 10.74039626121521, $null),
    @(55, 2209, 0.9880876541137695, $null),
    @(56, 6919, 5.808436393737793, $null),
    @(57, 1646, 0.7998125553131104, $null),
    @(58, 1441, 0.6127204895019531, $null),
    @(59, 2388, 1.172754287719727, $null),
    @(60, 5613, 3.470796346664429, $null),
    @(61, 21590, 43.87055015563965, $null),
    @(62, 10808, 12.96010136604309, $null),
    @(63, 3368, 1.703220367431641, $null),
    @(64, 849, 0.3763980865478516, 1),
    @(65, 3015, 1.754358530044556, $null),
    @(66, 7089, 6.641978740692139, $null),
    @(67, 8476, 9.046791315078735, $null),
    @(68, 8530, 7.968752145767212, $null),
    @(69, 6029, 4.761534690856934, 1),
    @(70, 2194, 0.6672215461730957, $null),
    @(71, 10483, 8.42224645614624, $null),
    @(72, 2198, 1.099310636520386, $null),
    @(73, 2304, 0.8563399314880371, $null),
    @(74, 4371, 2.850404024124146, $null),
    @(75, 9433, 10.06337308883667, $null),
    @(76, 6886, 5.758676767349243, $null),
    @(77, 9631, 8.981184482574463, $null),
    @(78, 9986, 10.87247657775879, $null),
    @(79, 1861, 0.8603420257568359, $null),
    @(80, 3642, 2.224096298217773, $null),
    @(81, 4908, 3.357852935791016, 1),
    @(82, 13979, 17.2071430683136, $null),
    @(83, 4865, 3.259887218475342, $null),
    @(84, 1923, 0.9256594181060791, $null),
    @(85, 5418, 3.977680683135986, $null),
    @(86, 1822, 0.8170063495635986, $null),
    @(87, 16770, 25.4084038734436, $null),
    @(88, 2454, 1.091157913208008, $null),
    @(89, 12218, 14.84336280822754, $null),
    @(90, 8483, 7.824760437011719, $null),
    @(91, 7254, 6.483772993087769, 1),
    @(92, 14943, 20.92472052574158, $null),
    @(93, 1851, 0.8043766021728516, $null),
    @(94, 8277, 7.748855829238892, $null),
    @(95, 4441, 2.877262353897095, $null),
    @(96, 3901, 2.430032253265381, $null),
    @(97, 1884, 0.8806445598602295, $null),
    @(98, 2565, 1.303118944168091, $null),
    @(99, 1558, 0.606403112411499, $null),
    @(100, 9021, 7.516368865966797, $null),
    @(101, 134, 0.03957748413085938, $null)
)

foreach ($u in $updates) {
    $r = $u[0]
    $iterations = $u[1]
    $runtime = $u[2]
    $success = $u[3]

    $ws.Cells.Item($r, 2).Value = $iterations
    $ws.Cells.Item($r, 3).Value = $runtime
    if ($success -ne $null) {
        $ws.Cells.Item($r, 4).Value = $success
    }
}

# Update Summary sheet row 2: A2 Iterations Median, B2 Runtime Median, C2 Success Ratio
$wsSummary.Range("A2").Value = 5163
$wsSummary.Range("B2").Value = 3.838541626930237
$wsSummary.Range("C2").Value = 0.99

Write-Host "Update complete"
